# Add the I0 (column I) and IF (column J) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, styled like the other header cells
# (copy the formatting from the existing H1 header cell so the new
# headers share the same bold/bordered/centered style).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-48: row number, I value, J value
$data = @(
    @(2, 7, 7),
    @(3, 7, 7),
    @(4, 6, 6),
    @(5, 9, 9),
    @(6, 8, 8),
    @(7, 7, 7),
    @(8, 7, 7),
    @(9, 7, 7),
    @(10, 7, 7),
    @(11, 7, 7),
    @(12, 8, 8),
    @(13, 9, 9),
    @(14, 7, 7),
    @(15, 7, 7),
    @(16, 7, 7),
    @(17, 8, 8),
    @(18, 8, 8),
    @(19, 7, 7),
    @(20, 9, 9),
    @(21, 9, 9),
    @(22, 8, 8),
    @(23, 8, 8),
    @(24, 8, 8),
    @(25, 7, 7),
    @(26, 9, 9),
    @(27, 10, 10),
    @(28, 9, 9),
    @(29, 10, 10),
    @(30, 8, 8),
    @(31, 8, 8),
    @(32, 9, 9),
    @(33, 7, 8),
    @(34, 8, 8),
    @(35, 8, 8),
    @(36, 7, 7),
    @(37, 8, 8),
    @(38, 7, 7),
    @(39, 8, 8),
    @(40, 8, 8),
    @(41, 2, 2),
    @(42, 6, 6),
    @(43, 5, 5),
    @(44, 9, 9),
    @(45, 6, 6),
    @(46, 6, 6),
    @(47, 6, 6),
    @(48, 5, 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
